# Adds a "Creator" (Virginia Scarlett) identifier/affiliation block to the
# csv2json attribute-table sheet: nameIdentifiers / nameIdentifierScheme /
# schemeURI / two Affiliations rows, inserted right after her existing
# "nameType" row (so they land before the second creator, William
# Shakespeare). Also re-applies the font formatting that Excel picked up
# when this data was pasted in from its source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: rows 4-8 (5 new rows) are inserted, pushing the old rows
# 4-8 down to 9-13.
[void]$ws.Rows("4:8").Insert()

# --- New rows 4-8: all belong to creator #1 (Virginia Scarlett) ---
# Column A ("creators") and column B (creator index 1) first.
$ws.Range("A4").Value = "creators"
$ws.Range("B4").Value = 1
$ws.Range("A5").Value = "creators"
$ws.Range("B5").Value = 1
$ws.Range("A6").Value = "creators"
$ws.Range("B6").Value = 1
$ws.Range("A7").Value = "creators"
$ws.Range("B7").Value = 1
$ws.Range("A8").Value = "creators"
$ws.Range("B8").Value = 1

# Column C: the Attr_key for each new row.
$ws.Range("C4").Value = "nameIdentifiers"
$ws.Range("C5").Value = "nameIdentifierScheme"
$ws.Range("C6").Value = "schemeURI"
$ws.Range("C7").Value = "Affiliations"
$ws.Range("C8").Value = "Affiliations"

# Column D: the Attr_value for each new row.
$ws.Range("D4").Value = "0000-0002-4156-2849"
$ws.Range("D5").Value = "ORCID"
$ws.Range("D6").Value = "https://orcid.org"
$ws.Range("D7").Value = "University of California, Berkeley"
$ws.Range("D8").Value = "HHMI Janelia Research Campus"

# --- Re-apply font formatting across the whole table (A1:D13) ---
# Column A throughout uses an 11pt black font.
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 0
$ws.Range("A1:A13").Font.Size = 11
$ws.Range("A1:A13").Font.Color = 0

# B:D for rows 1-10 use an 11pt dark-grey (#24292F) font.
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").Font.Color = 3090724
$ws.Range("B1:D10").Font.Size = 11
$ws.Range("B1:D10").Font.Color = 3090724

# B:D for the trailing rows (11-13: titles/publisher/publicationYear) keep
# the plain 11pt theme-colored font.
$ws.Range("B11").Font.Size = 11
$ws.Range("B11:D13").Font.Size = 11

# Match the saved selection cursor.
[void]$ws.Range("G9").Select()

Write-Host "done"
